$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35; this shifts the existing rows 35..63
# down to 36..64, extending the used range from A1:R63 to A1:R64.
$ws.Rows("35:35").Insert()

# Populate the newly inserted row 35 with the new Espárragos record.
$ws.Range("A35").Value = 5
$ws.Range("B35").Value = 'Macroferia Regional de Talca'
$ws.Range("C35").Value = 'Maule'
$ws.Range("D35").Value2 = 44827
$ws.Range("E35").Value = 7
$ws.Range("F35").Value = 300000000
$ws.Range("G35").Value = 'Espárragos'
$ws.Range("H35").Value = 'Sin especificar'
$ws.Range("I35").Value = 'Primera'
$ws.Range("J35").Value = 1000
$ws.Range("K35").Value = 2000
$ws.Range("L35").Value = 2000
$ws.Range("M35").Value = 2000
$ws.Range("N35").Value = '$/kilo'
$ws.Range("O35").Value = 'Provincia de Linares'
$ws.Range("P35").Value = 2000
$ws.Range("Q35").Value = 1
$ws.Range("R35").Value = 'Hortaliza'
